# Apply the September 28, 2021 LEGO search-results row reshuffle.
# Several newly-scraped rows ("Hogwarts(TM) Moment: Potions Class",
# "Hogwarts(TM) Wizard's Chess", "Hogwarts(TM): Polyjuice Potion Mistake",
# "Hogwarts(TM): Fluffy Encounter", "Hogsmeade(TM) Village Visit",
# "The Rise of Voldemort(TM)", "Harry Potter(TM) & Hedwig(TM)", and
# "Harry Potter(TM) and Fantastic Beasts(TM)") were interleaved into the
# existing list, which shifts the Name/Price/Availability values of the
# affected rows. Only the cells whose value actually changes are touched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds "$xx.xx"-style text; pre-mark it as Text so the
# COM layer does not coerce the new strings into numbers.
foreach ($addr in @("B11","B12","B13","B17","B18","B19","B21","B22","B23","B24","B25","B26","B30","B31","B32","B40")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Cells.Item(11, 1).Value = 'Hogwarts™ Whomping Willow™'
$ws.Cells.Item(11, 2).Value = '$69.99'
$ws.Cells.Item(12, 1).Value = 'Hogwarts™ Room of Requirement'
$ws.Cells.Item(12, 2).Value = '$19.99'
$ws.Cells.Item(13, 1).Value = 'Hagrid''s Hut: Buckbeak''s Rescue'
$ws.Cells.Item(13, 2).Value = '$59.99'
$ws.Cells.Item(17, 1).Value = 'Hogwarts™ Moment: Potions Class'
$ws.Cells.Item(17, 2).Value = '$29.99'
$ws.Cells.Item(18, 1).Value = 'Hogwarts™ Students Acc. Set'
$ws.Cells.Item(18, 2).Value = '$14.99'
$ws.Cells.Item(19, 1).Value = 'Voldemort™, Nagini & Bellatrix'
$ws.Cells.Item(19, 2).Value = '$24.99'
$ws.Cells.Item(21, 1).Value = 'Harry Potter™ Key Chain'
$ws.Cells.Item(21, 2).Value = '$5.99'
$ws.Cells.Item(22, 1).Value = 'Harry Potter™ Hogwarts™ Crests'
$ws.Cells.Item(22, 2).Value = '$119.99'
$ws.Cells.Item(23, 1).Value = 'Hogwarts™: Polyjuice Potion Mistake'
$ws.Cells.Item(23, 2).Value = '$19.99'
$ws.Cells.Item(24, 1).Value = 'Hogwarts™ Icons - Collectors'' Edition'
$ws.Cells.Item(24, 2).Value = '$249.99'
$ws.Cells.Item(25, 1).Value = 'The Knight Bus™'
$ws.Cells.Item(25, 2).Value = '$39.99'
$ws.Cells.Item(26, 1).Value = 'Hogwarts™ Moment: Charms Class'
$ws.Cells.Item(26, 2).Value = '$29.99'
$ws.Cells.Item(27, 1).Value = 'Hogwarts™: First Flying Lesson'
$ws.Cells.Item(28, 1).Value = 'Hogwarts™ Moment: Transfiguration Class'
$ws.Cells.Item(30, 1).Value = 'Hogwarts™: Fluffy Encounter'
$ws.Cells.Item(30, 2).Value = '$39.99'
$ws.Cells.Item(31, 1).Value = 'Hogsmeade™ Village Visit'
$ws.Cells.Item(31, 2).Value = '$79.99'
$ws.Cells.Item(32, 1).Value = 'Diagon Alley™'
$ws.Cells.Item(32, 2).Value = '$399.99'
$ws.Cells.Item(37, 1).Value = 'The Rise of Voldemort™'
$ws.Cells.Item(37, 2).Value = 'N/A'
$ws.Cells.Item(37, 3).Value = 'Retired Product'
$ws.Cells.Item(38, 1).Value = 'Harry Potter™ & Hedwig™'
$ws.Cells.Item(38, 2).Value = 'N/A'
$ws.Cells.Item(38, 3).Value = 'Retired Product'
$ws.Cells.Item(40, 1).Value = 'Hogwarts™ Great Hall'
$ws.Cells.Item(40, 2).Value = '$99.99'
$ws.Cells.Item(40, 3).Value = 'Sold out'
